# Run-mode parameterization: insert a new "TestSuite" sheet at the front of
# the workbook that drives which test cases (AddCustomer / OpenAccount) run,
# plus a global "loginAsManager" toggle.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts immediately before the active sheet, so adding it
# now (while AddCustomer - the first/active sheet - is still selected) puts
# it first, matching the target tab order: TestSuite, AddCustomer, OpenAccount.
$testSuite = $wb.Worksheets.Add()
$testSuite.Name = "TestSuite"

# Header row.
$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "RunMode"

# Per-test rows first (so the shared-string pool fills AddCustomer/OpenAccount/y/n
# ahead of the later-inserted loginAsManager row).
$testSuite.Range("A3").Value = "AddCustomer"
$testSuite.Range("A4").Value = "OpenAccount"
$testSuite.Range("B3").Value = "y"
$testSuite.Range("B4").Value = "n"

# Global login-mode row, inserted last (row 2, between header and test rows).
$testSuite.Range("A2").Value = "loginAsManager"
$testSuite.Range("B2").Value = "y"

# Column widths to roughly match the authored layout.
$testSuite.Columns.Item(1).ColumnWidth = 18.5
$testSuite.Columns.Item(2).ColumnWidth = 25.35

# Leave the selection on A2, like the source workbook.
$testSuite.Range("A2").Select() | Out-Null
